# Re-sort the worksheet tabs: "总计" (summary) should come before "2021-Q2" (detail).
$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2 = $wb.Worksheets.Item("2021-Q2")

# Move the "总计" sheet so it becomes the first tab (ahead of "2021-Q2").
$total.Move($q2)
